$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price column to Text so numeric-looking values
# (e.g. "100.00") are stored as literal strings instead of being
# auto-converted to numbers by Excel's type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '50.927.35'
$ws.Range("E2").Value = '  -1.01%  '
$ws.Range("D3").Value = '2.926.89'
$ws.Range("E3").Value = '  -1.38%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '372.92'
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").Value = '100.00'
$ws.Range("E6").Value = '  -4.41%  '
$ws.Range("D7").Value = '0.531'
$ws.Range("E7").Value = '  -1.79%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.576'
$ws.Range("E9").Value = '  -2.44%  '
$ws.Range("D10").Value = '35.82'
$ws.Range("E10").Value = '  -3.47%  '
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").Value = '0.0839'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").Value = '3.387.11'
$ws.Range("E13").Value = '  -1.55%  '
$ws.Range("D14").Value = '17.86'
$ws.Range("E14").Value = '  -2.63%  '
$ws.Range("D15").Value = '7.40'
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = '11.19'
$ws.Range("E16").Value = '  +51.50%  '
$ws.Range("D17").Value = '2.918.34'
$ws.Range("E17").Value = '  -1.78%  '
$ws.Range("D18").Value = '0.959'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("D19").Value = '50.887.45'
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").Value = '3.13'
$ws.Range("E20").Value = '  -5.75%  '
$ws.Range("D21").Value = '12.22'
$ws.Range("E21").Value = '  -5.20%  '
$ws.Range("D22").Value = '0.0₃0948'
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("D23").Value = '263.08'
$ws.Range("E23").Value = '  +0.67%  '
$ws.Range("D24").Value = '68.18'
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("D25").Value = '3.14'
$ws.Range("E25").Value = '  +10.74%  '
$ws.Range("D26").Value = '8.01'
$ws.Range("E26").Value = '  -1.68%  '
$ws.Range("D27").Value = '7.20'
$ws.Range("E27").Value = '  -4.30%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = '0.110'
$ws.Range("E29").Value = '  -2.78%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '25.38'
$ws.Range("E30").Value = '  -1.75%  '
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = '0.162'
$ws.Range("E31").Value = '  -5.06%  '
$ws.Range("D32").Value = '9.82'
$ws.Range("E32").Value = '  -0.55%  '
$ws.Range("D33").Value = '50.47'
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("E34").Value = '  -2.82%  '
$ws.Range("D35").Value = '32.74'
$ws.Range("E35").Value = '  -6.26%  '
$ws.Range("E36").Value = '  -1.99%  '
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").Value = '3.11'
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("E39").Value = '  -0.90%  '
$ws.Range("D40").Value = '16.14'
$ws.Range("E40").Value = '  -5.63%  '
$ws.Range("D41").Value = '1.77'
$ws.Range("E41").Value = '  -4.17%  '
$ws.Range("E42").Value = '  -5.70%  '
$ws.Range("E43").Value = '  -4.60%  '
$ws.Range("D44").Value = '20.89'
$ws.Range("E44").Value = '  -4.10%  '
$ws.Range("E45").Value = '  -1.67%  '
$ws.Range("D46").Value = '0.272'
$ws.Range("E46").Value = '  -5.92%  '
$ws.Range("D47").Value = '3.22'
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("D49").Value = '1.972.27'
$ws.Range("E49").Value = '  -2.95%  '
$ws.Range("D50").Value = '0.0323'
$ws.Range("E50").Value = '  -5.18%  '
$ws.Range("D51").Value = '1.30'
$ws.Range("E51").Value = '  +1.43%  '

# Restore the original (default/General) formatting on the Price column
# now that the text values are safely stored.
$ws.Range("D2:D51").ClearFormats()
